$wb = $excel.ActiveWorkbook

# --- Input sheet: insert a new "repeat every" row before the existing
#     "selectweekdaysfriday" row, and restate the weekday value as "FRI" ---
$wsInput = $wb.Worksheets.Item("Input")

# Push the "selectweekdaysfriday" row (and everything after it) down one
# row, opening up row 11 for the new "selectRepeatevery" field.
$wsInput.Rows(11).Insert()

# New row 11: selectRepeatevery = 1
$wsInput.Range("A11").Value = "selectRepeatevery"
$wsInput.Range("B11").Value = 1

# Row 12 (previously row 11): selectweekdaysfriday now stores "FRI" and
# loses its green highlight fill (back to the default/no style).
$wsInput.Range("B12").Value = "FRI"
$wsInput.Range("B12").Style = "Normal"

# Row 10's "weekly" value also loses its green highlight fill.
$wsInput.Range("B10").Style = "Normal"

# Selection on the Input sheet moves to the newly added cell.
$wsInput.Range("B11").Select()

# --- Repayment schedule sheet: update the remembered selection ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("P13").Select()

# --- Summary sheet becomes the active tab with a new selection ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("F3").Select()
